# Update the Diebold-Mariano results table ("analisis de las 3 primeras
# simulaciones"): refresh DM_stat (col C) and p_valor (col D) for every
# data row (2-73), and flip row 70's "Significativo" (col E) from "No" to
# "Si" now that its recomputed p-value is below the 0.05 threshold.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"5.636518091129741"
$ws.Range("D2").Value = [double]"1.145354048293079E-05"
$ws.Range("C3").Value = [double]"4.311866060786169"
$ws.Range("D3").Value = [double]"0.0002816624389738731"
$ws.Range("C4").Value = [double]"5.791045711614122"
$ws.Range("D4").Value = [double]"7.951149563867332E-06"
$ws.Range("C5").Value = [double]"6.120908812012834"
$ws.Range("D5").Value = [double]"3.67921254418313E-06"
$ws.Range("C6").Value = [double]"6.101716454321259"
$ws.Range("D6").Value = [double]"3.84664969921289E-06"
$ws.Range("C7").Value = [double]"7.91837218254468"
$ws.Range("D7").Value = [double]"6.995991586045136E-08"
$ws.Range("C8").Value = [double]"5.793924396380716"
$ws.Range("D8").Value = [double]"7.897454239236268E-06"
$ws.Range("C9").Value = [double]"7.222888882224355"
$ws.Range("D9").Value = [double]"3.083063053743729E-07"
$ws.Range("C10").Value = [double]"3.979940955020395"
$ws.Range("D10").Value = [double]"0.0006333609825530662"
$ws.Range("C11").Value = [double]"4.505959219671162"
$ws.Range("D11").Value = [double]"0.0001753018583174892"
$ws.Range("C12").Value = [double]"5.406969511476652"
$ws.Range("D12").Value = [double]"1.97813312365902E-05"
$ws.Range("C13").Value = [double]"6.255841783842101"
$ws.Range("D13").Value = [double]"2.69398390417841E-06"
$ws.Range("C14").Value = [double]"5.99361615974306"
$ws.Range("D14").Value = [double]"4.946358789670313E-06"
$ws.Range("C15").Value = [double]"8.04942311779922"
$ws.Range("D15").Value = [double]"5.328705965368385E-08"
$ws.Range("C16").Value = [double]"5.713297110069927"
$ws.Range("D16").Value = [double]"9.551046070699343E-06"
$ws.Range("C17").Value = [double]"7.272712060199441"
$ws.Range("D17").Value = [double]"2.766339517457084E-07"
$ws.Range("C18").Value = [double]"5.647896610076669"
$ws.Range("D18").Value = [double]"1.114891323972245E-05"
$ws.Range("C19").Value = [double]"4.299167632376165"
$ws.Range("D19").Value = [double]"0.0002905397748014238"
$ws.Range("C20").Value = [double]"5.817629428354874"
$ws.Range("D20").Value = [double]"7.46910221427477E-06"
$ws.Range("C21").Value = [double]"6.12799034992442"
$ws.Range("D21").Value = [double]"3.619328720061432E-06"
$ws.Range("C22").Value = [double]"6.114020870926963"
$ws.Range("D22").Value = [double]"3.738430381039493E-06"
$ws.Range("C23").Value = [double]"7.909158998473044"
$ws.Range("D23").Value = [double]"7.13178889277799E-08"
$ws.Range("C24").Value = [double]"5.799345440511664"
$ws.Range("D24").Value = [double]"7.797337738812971E-06"
$ws.Range("C25").Value = [double]"7.217805245870291"
$ws.Range("D25").Value = [double]"3.11740962821716E-07"
$ws.Range("C26").Value = [double]"2.494871527303883"
$ws.Range("D26").Value = [double]"0.02059645052341619"
$ws.Range("C27").Value = [double]"4.643670539851906"
$ws.Range("D27").Value = [double]"0.0001252677480902964"
$ws.Range("C28").Value = [double]"5.675215368924245"
$ws.Range("D28").Value = [double]"1.045074495409359E-05"
$ws.Range("C29").Value = [double]"5.720821364936925"
$ws.Range("D29").Value = [double]"9.38282756179909E-06"
$ws.Range("C30").Value = [double]"6.047728025981449"
$ws.Range("D30").Value = [double]"4.360638136802208E-06"
$ws.Range("C31").Value = [double]"7.706918887139473"
$ws.Range("D31").Value = [double]"1.090728460884094E-07"
$ws.Range("C32").Value = [double]"5.799533888521103"
$ws.Range("D32").Value = [double]"7.79388080918153E-06"
$ws.Range("C33").Value = [double]"7.203584202660262"
$ws.Range("D33").Value = [double]"3.215595709704644E-07"
$ws.Range("C34").Value = [double]"6.533099618323958"
$ws.Range("D34").Value = [double]"1.429661666385584E-06"
$ws.Range("C35").Value = [double]"3.139397787917392"
$ws.Range("D35").Value = [double]"0.004764677455612487"
$ws.Range("C36").Value = [double]"6.144802766344924"
$ws.Range("D36").Value = [double]"3.481112960557198E-06"
$ws.Range("C37").Value = [double]"5.888543658442072"
$ws.Range("D37").Value = [double]"6.323707134603396E-06"
$ws.Range("C38").Value = [double]"5.833610302955822"
$ws.Range("D38").Value = [double]"7.19375618696283E-06"
$ws.Range("C39").Value = [double]"8.127830936110364"
$ws.Range("D39").Value = [double]"4.532758457465036E-08"
$ws.Range("C40").Value = [double]"4.959639417898881"
$ws.Range("D40").Value = [double]"5.809569839931861E-05"
$ws.Range("C41").Value = [double]"9.549830277463794"
$ws.Range("D41").Value = [double]"2.772182927657241E-09"
$ws.Range("C42").Value = [double]"9.32770205465871"
$ws.Range("D42").Value = [double]"4.216664573064577E-09"
$ws.Range("C43").Value = [double]"4.462958617510705"
$ws.Range("D43").Value = [double]"0.0001947137330660098"
$ws.Range("C44").Value = [double]"5.745753801445408"
$ws.Range("D44").Value = [double]"8.846669779405403E-06"
$ws.Range("C45").Value = [double]"6.04567001826597"
$ws.Range("D45").Value = [double]"4.38156400450751E-06"
$ws.Range("C46").Value = [double]"5.858055782967347"
$ws.Range("D46").Value = [double]"6.792433822955246E-06"
$ws.Range("C47").Value = [double]"7.965213233173771"
$ws.Range("D47").Value = [double]"6.345664749218827E-08"
$ws.Range("C48").Value = [double]"5.706611632530271"
$ws.Range("D48").Value = [double]"9.703087466617077E-06"
$ws.Range("C49").Value = [double]"7.165013066362093"
$ws.Range("D49").Value = [double]"3.498223013487234E-07"
$ws.Range("C50").Value = [double]"6.344614883164578"
$ws.Range("D50").Value = [double]"2.197108210477694E-06"
$ws.Range("C51").Value = [double]"4.456080084220476"
$ws.Range("D51").Value = [double]"0.0001980130164593863"
$ws.Range("C52").Value = [double]"5.620627292294299"
$ws.Range("D52").Value = [double]"1.189320252614756E-05"
$ws.Range("C53").Value = [double]"6.025625107413418"
$ws.Range("D53").Value = [double]"4.590822228811575E-06"
$ws.Range("C54").Value = [double]"6.005243624812292"
$ws.Range("D54").Value = [double]"4.814066743019652E-06"
$ws.Range("C55").Value = [double]"8.012173778048647"
$ws.Range("D55").Value = [double]"5.756053234051706E-08"
$ws.Range("C56").Value = [double]"5.693267506750887"
$ws.Range("D56").Value = [double]"1.001396999433801E-05"
$ws.Range("C57").Value = [double]"7.305415274211782"
$ws.Range("D57").Value = [double]"2.57681175064306E-07"
$ws.Range("C58").Value = [double]"4.233609761091651"
$ws.Range("D58").Value = [double]"0.0003410168527335777"
$ws.Range("C59").Value = [double]"4.286107696026305"
$ws.Range("D59").Value = [double]"0.0002999616353533963"
$ws.Range("C60").Value = [double]"5.666122478986837"
$ws.Range("D60").Value = [double]"1.067804402676842E-05"
$ws.Range("C61").Value = [double]"6.046966344571899"
$ws.Range("D61").Value = [double]"4.368371011587158E-06"
$ws.Range("C62").Value = [double]"6.011397089147072"
$ws.Range("D62").Value = [double]"4.745522187876716E-06"
$ws.Range("C63").Value = [double]"8.007952831548478"
$ws.Range("D63").Value = [double]"5.806658842644197E-08"
$ws.Range("C64").Value = [double]"5.746966964764621"
$ws.Range("D64").Value = [double]"8.821391913160781E-06"
$ws.Range("C65").Value = [double]"7.25010344191222"
$ws.Range("D65").Value = [double]"2.905696052657447E-07"
$ws.Range("C66").Value = [double]"0.4600762888447206"
$ws.Range("D66").Value = [double]"0.6499762046036415"
$ws.Range("C67").Value = [double]"0.8043674110064734"
$ws.Range("D67").Value = [double]"0.4297917795193822"
$ws.Range("C68").Value = [double]"-0.01821346673869894"
$ws.Range("D68").Value = [double]"0.9856327348448009"
$ws.Range("C69").Value = [double]"0.1206867793560985"
$ws.Range("D69").Value = [double]"0.9050345930354642"
$ws.Range("C70").Value = [double]"2.266184405537039"
$ws.Range("D70").Value = [double]"0.03362187756545043"
$ws.Range("E70").Value = "Sí"
$ws.Range("C71").Value = [double]"3.918077605372396"
$ws.Range("D71").Value = [double]"0.0007363391154457943"
$ws.Range("C72").Value = [double]"4.183102404599025"
$ws.Range("D72").Value = [double]"0.0003857997176008876"
$ws.Range("C73").Value = [double]"6.230586074180454"
$ws.Range("D73").Value = [double]"2.855354205744121E-06"

Write-Output "Applied 145 cell updates"
